$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Column F label: the shared string "12in\n" becomes "12in DW\n" (the
# trailing line break is preserved), reused by rows 2-5. A brand-new
# label "10in\n" is introduced for row 6.
# ---------------------------------------------------------------------
$ws.Range("F2:F4").Value = "12in DW`n"

# ---------------------------------------------------------------------
# Row 2
# ---------------------------------------------------------------------
$ws.Range("B2").Value = 1315.51
$ws.Range("C2").Value = 1312.09
$ws.Range("D2").Value = 3.42
$ws.Range("E2").Value = 3.42

# ---------------------------------------------------------------------
# Row 3
# ---------------------------------------------------------------------
$ws.Range("B3").Value = 1318.29
$ws.Range("C3").Value = 1314.91
$ws.Range("D3").Value = 3.37
$ws.Range("E3").Value = 3.38

# ---------------------------------------------------------------------
# Row 4 - E4 also switches formatting from the bold/red "mismatch" style
# to the plain numeric style (it now matches D4), so copy D4's format
# across after updating the value.
# ---------------------------------------------------------------------
$ws.Range("B4").Value = 1318.7
$ws.Range("C4").Value = 1315.15
$ws.Range("D4").Value = 3.55
$ws.Range("E4").Value = 3.55
$ws.Range("D4").Copy()
$ws.Range("E4").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# New row 5 - clone row 4's formatting, then fill in the new values.
# ---------------------------------------------------------------------
$ws.Range("A4:F4").Copy()
$ws.Range("A5:F5").PasteSpecial(-4122)
$ws.Range("A5").Value = 400
$ws.Range("B5").Value = 1320.92
$ws.Range("C5").Value = 1315.96
$ws.Range("D5").Value = 4.96
$ws.Range("E5").Value = 4.96
$ws.Range("F5").Value = "12in DW`n"

# ---------------------------------------------------------------------
# New row 6 - same cloned formatting, with the new "10in" tile label.
# ---------------------------------------------------------------------
$ws.Range("A4:F4").Copy()
$ws.Range("A6:F6").PasteSpecial(-4122)
$ws.Range("A6").Value = 500
$ws.Range("B6").Value = 1324.09
$ws.Range("C6").Value = 1316.06
$ws.Range("D6").Value = 8.03
$ws.Range("E6").Value = 8.03
$ws.Range("F6").Value = "10in`n"

$excel.CutCopyMode = $false

# Setting the multi-line "DW"/"10in" labels makes the host auto-grow the
# row height; AutoFit puts every touched row back to the sheet's
# (non-custom) standard height so the rows serialize the same as before.
$ws.Rows("2:6").AutoFit()

